# Apply updated GHI/DNI/DHI simulation values for the RAAL Production model.
$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

# --- Daily sheet, row 2 ---
$daily.Range("G2").Value = 2652.19
$daily.Range("H2").Value = 5890.93
$daily.Range("I2").Value = 686.78
$daily.Range("J2").Value = 665.4400000000001
$daily.Range("L2").Value = 665.4400000000001

# --- Hourly sheet ---
# row 9
$hourly.Range("I9").Value = 22.61

# row 10
$hourly.Range("I10").Value = 377.56
$hourly.Range("K10").Value = 21.28
$hourly.Range("M10").Value = 21.28

# row 11
$hourly.Range("H11").Value = 219.89
$hourly.Range("I11").Value = 606.79
$hourly.Range("K11").Value = 55.13
$hourly.Range("M11").Value = 55.13

# row 12
$hourly.Range("H12").Value = 337.72
$hourly.Range("I12").Value = 712.29
$hourly.Range("J12").Value = 82.5
$hourly.Range("K12").Value = 84.56
$hourly.Range("M12").Value = 84.56

# row 13
$hourly.Range("I13").Value = 763.86
$hourly.Range("K13").Value = 104.38
$hourly.Range("M13").Value = 104.38

# row 14
$hourly.Range("H14").Value = 447.85
$hourly.Range("I14").Value = 781.2
$hourly.Range("J14").Value = 92.95999999999999
$hourly.Range("K14").Value = 112.26
$hourly.Range("M14").Value = 112.26

# row 15
$hourly.Range("H15").Value = 426.56
$hourly.Range("I15").Value = 769.47
$hourly.Range("K15").Value = 107.45
$hourly.Range("M15").Value = 107.45

# row 16
$hourly.Range("H16").Value = 355.59
$hourly.Range("I16").Value = 725.12
$hourly.Range("J16").Value = 84.3
$hourly.Range("K16").Value = 89.54000000000001
$hourly.Range("M16").Value = 89.54000000000001

# row 17
$hourly.Range("H17").Value = 243.71
$hourly.Range("I17").Value = 632.27
$hourly.Range("K17").Value = 60.93
$hourly.Range("M17").Value = 60.93

# row 18
$hourly.Range("H18").Value = 108.79
$hourly.Range("I18").Value = 437.09

# row 19
$hourly.Range("H19").Value = 8.300000000000001
$hourly.Range("I19").Value = 62.66

$wb.Save()
